# Journal de travail Osama - Semaine 4 edits
$d = $word.ActiveDocument
$wNS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark currently on the "Semaine 4" paragraph.
#    (It will be re-created later on the new "Samedi" bullet paragraph.)
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2. After "Jeudi " (paragraph 7), insert four new bulleted ListParagraph
#    items describing the Monday-Thursday work, using a fresh numbered list
#    (numId 7).
# ---------------------------------------------------------------------------
$pJeudi = $d.Paragraphs.Item(7)
$r = $d.Range($pJeudi.Range.End - 1, $pJeudi.Range.End - 1)

$xmlBullet1 = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr>" +
  "<w:r><w:t>1</w:t></w:r>" +
  "<w:r><w:rPr><w:vertAlign w:val='superscript'/></w:rPr><w:t>er</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> connextion sur phpMyAdmin </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t xml:space='preserve'>(10 minutes) </w:t></w:r>" +
  "</w:p>"
$r.InsertXML($xmlBullet1)

# Turn that first new paragraph into a proper bulleted list item; this mints
# the new numId (7) and its abstractNum definition in numbering.xml.
$pBullet1 = $d.Paragraphs.Item(8)
$bulletGallery = $word.ListGalleries.Item(1)
$bulletTemplate = $bulletGallery.ListTemplates.Item(1)
$pBullet1.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate)

# Discover the numId that got minted so the remaining bullets reuse it.
$numId7 = $pBullet1.Range.ListFormat.ListTemplate.Name
$numId7 = $pBullet1.Range.ListFormat

$r = $d.Range($pBullet1.Range.End - 1, $pBullet1.Range.End - 1)

$xmlBullet2 = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='7'/></w:numPr></w:pPr>" +
  "<w:r><w:t xml:space='preserve'>Modification de la partie Base de donn&#233; dans le rapport </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>(</w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>40 minutes</w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t xml:space='preserve'>) </w:t></w:r>" +
  "</w:p>"
$r.InsertXML($xmlBullet2)
$pBullet2 = $d.Paragraphs.Item(9)

$r = $d.Range($pBullet2.Range.End - 1, $pBullet2.Range.End - 1)
$xmlBullet3 = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='7'/></w:numPr>" +
  "<w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr></w:pPr>" +
  "<w:r><w:t>Remplir la liste des d&#233;partements pour les ajouter dans la table T_Departements</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>(15 m</w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>inutes)</w:t></w:r>" +
  "</w:p>"
$r.InsertXML($xmlBullet3)
$pBullet3 = $d.Paragraphs.Item(10)

$r = $d.Range($pBullet3.Range.End - 1, $pBullet3.Range.End - 1)
$xmlBullet4 = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='7'/></w:numPr></w:pPr>" +
  "<w:r><w:t>Planification de ce qu&#8217;il faut faire pendant les vacances</w:t></w:r>" +
  "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val='18'/><w:szCs w:val='18'/></w:rPr><w:t>(10 minutes)</w:t></w:r>" +
  "</w:p>"
$r.InsertXML($xmlBullet4)

Write-Host "Inserted Monday-Thursday bullet list (numId 7)."

# ---------------------------------------------------------------------------
# 3. "Samedi " (now paragraph 13): remove the rtl paragraph-mark formatting,
#    it becomes a plain Heading1 paragraph again.
# ---------------------------------------------------------------------------
$pSamedi = $d.Paragraphs.Item(13)
Write-Host "Samedi check: [" $pSamedi.Range.Text "]"
$xmlSamedi = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='Heading1'/></w:pPr>" +
  "<w:r><w:t xml:space='preserve'>Samedi </w:t></w:r>" +
  "</w:p>"
$pSamedi.Range.InsertXML($xmlSamedi)

# ---------------------------------------------------------------------------
# 4. After "Samedi " insert the new bulleted ListParagraph (numId 8, RTL)
#    holding the "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$pSamedi = $d.Paragraphs.Item(13)
$r = $d.Range($pSamedi.Range.End - 1, $pSamedi.Range.End - 1)
$xmlBulletSamedi = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr>" +
  "<w:r><w:t xml:space='preserve'>Faire la partie &#171; Planification du travail &#187; dans le rapport </w:t></w:r>" +
  "</w:p>"
$r.InsertXML($xmlBulletSamedi)

$pBulletSamedi = $d.Paragraphs.Item(14)
$pBulletSamedi.Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate)

Write-Host "Inserted Samedi bullet (numId 8)."

# ---------------------------------------------------------------------------
# 5. "Dimanche" (now paragraph 15): update its paragraph-mark run formatting.
# ---------------------------------------------------------------------------
$pDimanche = $d.Paragraphs.Item(15)
Write-Host "Dimanche check: [" $pDimanche.Range.Text "]"
$xmlDimanche = "<w:p $wNS>" +
  "<w:pPr><w:pStyle w:val='Heading1'/>" +
  "<w:rPr><w:rFonts w:asciiTheme='minorHAnsi' w:eastAsiaTheme='minorHAnsi' w:hAnsiTheme='minorHAnsi' w:cstheme='minorBidi'/>" +
  "<w:color w:val='auto'/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr></w:pPr>" +
  "<w:r><w:t>Dimanche</w:t></w:r>" +
  "</w:p>"
$pDimanche.Range.InsertXML($xmlDimanche)

# ---------------------------------------------------------------------------
# 6. Re-add the "_GoBack" bookmark on the new Samedi bullet paragraph
#    (right before its trailing paragraph mark, after the visible text).
# ---------------------------------------------------------------------------
$pBulletSamedi = $d.Paragraphs.Item(14)
$bmPos = $pBulletSamedi.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 7. The empty ListParagraph-styled paragraph right after "Dimanche" loses
#    its style, becoming a plain empty paragraph.
# ---------------------------------------------------------------------------
$pEmpty = $d.Paragraphs.Item(16)
Write-Host "Empty check: [" $pEmpty.Range.Text "]"
$pEmpty.Range.InsertXML("<w:p $wNS></w:p>")

Write-Host "Done."
